# Auto-generated data-driven edit script for chess game sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Alignment constants matching the existing move-cell style (style index 10:
# horizontal=left, vertical=center) so newly added cells match their siblings.
$xlHAlignLeft = -4131
$xlVAlignCenter = -4108

function Set-Move($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $text
    $cell.HorizontalAlignment = $xlHAlignLeft
    $cell.VerticalAlignment = $xlVAlignCenter
}

# --- Row 8: fix round number text ("2" -> "1") ---
$ws.Range("A8").Value = 'Играна во "1" туре 01.01.2000'
$ws.Range("G8").Value = 'Играна во "1" туре 01.01.2000'

# --- Rows 13-21: correct existing moves (cells already exist, keep their style) ---
$moveUpdates = @(
    @{Row=13; B='Кf3'; C='Кc6'; H='Кf3'; I='Кc6'},
    @{Row=14; B='d4'; C='f6'; H='c3'; I='d6'},
    @{Row=15; B='d5'; C='Кd4'; H='d4'; I='Cg4'},
    @{Row=16; B='Кxd4'; C='exd4'; H='Cg5'; I='f6'},
    @{Row=17; B='Фxd4'; C='d6'; H='Ce3'; I='Фd7'},
    @{Row=18; B='Кc3'; C='f5'; H='Кbd2'; I='O-O-O'},
    @{Row=19; B='e5'; C='dxe5'; H='h3'; I='Cxf3'},
    @{Row=20; B='Фxe5+'; C='Фe7'; H='Фxf3'; I='Кge7'},
    @{Row=21; B='Фxe7+'; C='Кxe7'; H='d5'; I='Фf7'}
)
foreach ($r in $moveUpdates) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
}

# --- Rows 22-43: add new moves (both white/black columns) ---
$moveAdditionsFull = @(
    @{Row=22; B='Cd3'; C='c6'; H='dxc6'; I='Кxc6'},
    @{Row=23; B='O-O'; C='Кxd5'; H='Cxa7'; I='Кxa7'},
    @{Row=24; B='Лe1+'; C='Ce7'; H='g4'; I='d5'},
    @{Row=25; B='Кxd5'; C='cxd5'; H='exd5'; I='Лxd5'},
    @{Row=26; B='h3'; C='Ce6'; H='Лd1'; I='Лc5'},
    @{Row=27; B='Лxe6'; C='O-O-O'; H='Фd3'; I='Cd6'},
    @{Row=28; B='Лxe7'; C='Лhe8'; H='Фf5+'; I='Крb8'},
    @{Row=29; B='Лxg7'; C='Лe1+'; H='b4'; I='Лxc4'},
    @{Row=30; B='Крh2'; C='h5'; H='Cxc4'; I='Фxc4'},
    @{Row=31; B='Cxf5+'; C='Крb8'; H='Крe2'; I='Фxa3'},
    @{Row=32; B='Лg5'; C='Лh8'; H='Фe4'; I='Фxb4'},
    @{Row=33; B='b3'; C='d4'; H='Фe3'; I='Фb5+'},
    @{Row=34; B='Cb2'; C='Лxa1'; H='Крf3'; I='Cc5'},
    @{Row=35; B='Cxa1'; C='d3'; H='Фa3'; I='Cxa3'},
    @{Row=36; B='cxd3'; C='h4'; H='Крg3'; I='Фc3+'},
    @{Row=37; B='d4'; C='Лf8'; H='f3'; I='Фxd2'},
    @{Row=38; B='d5'; C='Крc7'; H='f4'; I='Фe3+'},
    @{Row=39; B='Ce5+'; C='Крb6'; H='Крh4'; I='Фxf4'},
    @{Row=40; B='Cd4+'; C='Крc7'; H='Крh5'; I='g6+'},
    @{Row=41; B='Cg7'; C='Лf7'; H='Крh4'; I='g5+'},
    @{Row=42; B='g4'; C='Крd6'; H='Крh5'; I='Фg3'},
    @{Row=43; B='Ca1'; C='Крxd5'; H='Лd1'; I='Фxh3#'}
)
foreach ($r in $moveAdditionsFull) {
    Set-Move $r.Row 2 $r.B
    Set-Move $r.Row 3 $r.C
    Set-Move $r.Row 8 $r.H
    Set-Move $r.Row 9 $r.I
}

# --- Rows 44-47: add new moves (white side columns B/C only) ---
$moveAdditionsBC = @(
    @{Row=44; B='Лh5'; C='Лe7'},
    @{Row=45; B='Лxh4'; C='Лe2'},
    @{Row=46; B='Крg3'; C='Лxa2'},
    @{Row=47; B='Ch8'; C='Лa3'}
)
foreach ($r in $moveAdditionsBC) {
    Set-Move $r.Row 2 $r.B
    Set-Move $r.Row 3 $r.C
}

